# Auto-generated Excel COM-interop script to apply market data value updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1902.76
$ws.Range("I19").Value = 1523.25
$ws.Range("K19").Value = 1523.25
$ws.Range("M19").Value = -1348.25

$ws.Range("H62").Value = 2494
$ws.Range("I62").Value = 2494
$ws.Range("K62").Value = 2494
$ws.Range("M62").Value = -1870

$ws.Range("H65").Value = 2494
$ws.Range("I65").Value = 2494
$ws.Range("K65").Value = 12470
$ws.Range("M65").Value = -9350

$ws.Range("H80").Value = 1846.375
$ws.Range("I80").Value = 1686.8334
$ws.Range("J80").Value = 2325
$ws.Range("K80").Value = 5060.5002
$ws.Range("L80").Value = 6975
$ws.Range("M80").Value = -4062.5002
$ws.Range("N80").Value = -8971

$ws.Range("H83").Value = 1846.375
$ws.Range("I83").Value = 1686.8334
$ws.Range("J83").Value = 2325
$ws.Range("K83").Value = 15181.5006
$ws.Range("L83").Value = 20925
$ws.Range("M83").Value = -10189.5006
$ws.Range("N83").Value = -30909

$ws.Range("H106").Value = 37055260
$ws.Range("I106").Value = 40016440
$ws.Range("K106").Value = 40016440
$ws.Range("M106").Value = -40015809

$ws.Range("H116").Value = 4566
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4566
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4566
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11450

$ws.Range("H135").Value = 1974.091
$ws.Range("J135").Value = 3333
$ws.Range("L135").Value = 29997
$ws.Range("N135").Value = -35067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4747.222
$ws.Range("I32").Value = 3203.1035
$ws.Range("J32").Value = 11144.286
$ws.Range("K32").Value = 3203.1035
$ws.Range("L32").Value = 11144.286
$ws.Range("M32").Value = -2916.1035
$ws.Range("N32").Value = -11718.286

$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3288
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 3945.2727
$ws.Range("I74").Value = 3711.4443
$ws.Range("K74").Value = 3711.4443
$ws.Range("M74").Value = -2837.4443

$ws.Range("H77").Value = 3945.2727
$ws.Range("I77").Value = 3711.4443
$ws.Range("K77").Value = 18557.2215
$ws.Range("M77").Value = -14189.2215

$ws.Range("H132").Value = 3891.95
$ws.Range("I132").Value = 4028.6667
$ws.Range("J132").Value = 3686.875
$ws.Range("K132").Value = 12086.0001
$ws.Range("L132").Value = 11060.625
$ws.Range("M132").Value = -9556.000100000001
$ws.Range("N132").Value = -16120.625

$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 3500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7950
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4630.3
$ws.Range("I20").Value = 3247.75
$ws.Range("J20").Value = 5552
$ws.Range("K20").Value = 3247.75
$ws.Range("L20").Value = 5552
$ws.Range("M20").Value = -3000.75
$ws.Range("N20").Value = -6046

$ws.Range("H80").Value = 337.6154
$ws.Range("I80").Value = 341
$ws.Range("K80").Value = 341
$ws.Range("M80").Value = 657

$ws.Range("H83").Value = 337.6154
$ws.Range("I83").Value = 341
$ws.Range("K83").Value = 1705
$ws.Range("M83").Value = 3287

$ws.Range("H134").Value = 1804.4
$ws.Range("I134").Value = 1804.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5413.200000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2878.200000000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 91845.336
$ws.Range("J22").Value = 33467
$ws.Range("L22").Value = 33467
$ws.Range("N22").Value = -34167

$ws.Range("H31").Value = 3239.6453
$ws.Range("I31").Value = 2237.5652
$ws.Range("K31").Value = 2237.5652
$ws.Range("M31").Value = -1942.5652

$ws.Range("H34").Value = 3239.6453
$ws.Range("I34").Value = 2237.5652
$ws.Range("K34").Value = 2237.5652
$ws.Range("M34").Value = -2035.5652

$ws.Range("H99").Value = 13605.12
$ws.Range("I99").Value = 9948.799999999999
$ws.Range("J99").Value = 16042.667
$ws.Range("K99").Value = 9948.799999999999
$ws.Range("L99").Value = 16042.667
$ws.Range("M99").Value = -8450.799999999999
$ws.Range("N99").Value = -19038.667

$ws.Range("H122").Value = 1644.5
$ws.Range("I122").Value = 992.6
$ws.Range("J122").Value = 2296.4
$ws.Range("K122").Value = 2977.8
$ws.Range("L122").Value = 6889.200000000001
$ws.Range("M122").Value = -527.8000000000002
$ws.Range("N122").Value = -11789.2

$ws.Range("H126").Value = 13605.12
$ws.Range("I126").Value = 9948.799999999999
$ws.Range("J126").Value = 16042.667
$ws.Range("K126").Value = 29846.4
$ws.Range("L126").Value = 48128.001
$ws.Range("M126").Value = -27376.4
$ws.Range("N126").Value = -53068.001

$ws.Range("H132").Value = 2784.3809
$ws.Range("I132").Value = 2498.3125
$ws.Range("K132").Value = 7494.9375
$ws.Range("M132").Value = -4964.9375

$ws.Range("H134").Value = 2599.182
$ws.Range("I134").Value = 1718.8
$ws.Range("J134").Value = 3332.8333
$ws.Range("K134").Value = 5156.4
$ws.Range("L134").Value = 9998.499899999999
$ws.Range("M134").Value = -2621.4
$ws.Range("N134").Value = -15068.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 133.36842
$ws.Range("I12").Value = 166.875
$ws.Range("J12").Value = 109
$ws.Range("K12").Value = 500.625
$ws.Range("L12").Value = 327
$ws.Range("M12").Value = -327.625
$ws.Range("N12").Value = -673

$ws.Range("H86").Value = 371.2
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 371.2
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H113").Value = 299.66666
$ws.Range("I113").Value = 299.5
$ws.Range("K113").Value = 898.5
$ws.Range("M113").Value = 1271.5

$ws.Range("H122").Value = 333.2
$ws.Range("J122").Value = 359.66666
$ws.Range("L122").Value = 3236.99994
$ws.Range("N122").Value = -8136.99994

$ws.Range("H131").Value = 1410.1818
$ws.Range("I131").Value = 595.7143
$ws.Range("J131").Value = 1506.8136
$ws.Range("K131").Value = 1787.1429
$ws.Range("L131").Value = 4520.4408
$ws.Range("M131").Value = 3252.8571
$ws.Range("N131").Value = -14600.4408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7664.1665
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 7996.25
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 7996.25
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -8536.25

$ws.Range("H73").Value = 7664.1665
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 7996.25
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 7996.25
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -9868.25

$ws.Range("H113").Value = 27799422
$ws.Range("I113").Value = 50016360
$ws.Range("J113").Value = 28250
$ws.Range("K113").Value = 50016360
$ws.Range("L113").Value = 28250
$ws.Range("M113").Value = -50014190
$ws.Range("N113").Value = -32590

$ws.Range("H122").Value = 94084.45
$ws.Range("I122").Value = 2289.6
$ws.Range("K122").Value = 6868.799999999999
$ws.Range("M122").Value = -4418.799999999999

$ws.Range("H132").Value = 2995
$ws.Range("I132").Value = 2995
$ws.Range("K132").Value = 8985
$ws.Range("M132").Value = -6455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 567.26086
$ws.Range("I16").Value = 585.5714
$ws.Range("K16").Value = 585.5714
$ws.Range("M16").Value = -415.5714

$ws.Range("H82").Value = 113332.445
$ws.Range("J82").Value = 999999
$ws.Range("L82").Value = 999999
$ws.Range("N82").Value = -1000721

$ws.Range("H85").Value = 113332.445
$ws.Range("J85").Value = 999999
$ws.Range("L85").Value = 999999
$ws.Range("N85").Value = -1002495

$ws.Range("H136").Value = 5554.1333
$ws.Range("I136").Value = 5129.5454
$ws.Range("J136").Value = 6721.75
$ws.Range("K136").Value = 15388.6362
$ws.Range("L136").Value = 20165.25
$ws.Range("M136").Value = -12838.6362
$ws.Range("N136").Value = -25265.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2774.3704
$ws.Range("I132").Value = 2464.0454
$ws.Range("J132").Value = 4139.8
$ws.Range("K132").Value = 7392.1362
$ws.Range("L132").Value = 12419.4
$ws.Range("M132").Value = -4862.1362
$ws.Range("N132").Value = -17479.4
